$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.898.92'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '2.343.54'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.685'
$ws.Range("E5").Value = '  +3.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.48'
$ws.Range("E6").Value = '  +1.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.55'
$ws.Range("E7").Value = '  +3.28%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("E9").Value = '  +14.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  +1.81%  '

$ws.Range("E11").Value = '  +0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.90'
$ws.Range("E12").Value = '  +16.92%  '

$ws.Range("E13").Value = '  +11.02%  '

$ws.Range("E14").Value = '  +1.51%  '

$ws.Range("D15").Value = '2.695.93'
$ws.Range("E15").Value = '  -0.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.55'
$ws.Range("E16").Value = '  -1.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.915'
$ws.Range("E17").Value = '  +3.62%  '

$ws.Range("D18").Value = '2.348.56'
$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").Value = '43.785.78'
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("E20").Value = '  +0.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.62'
$ws.Range("E21").Value = '  +4.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.19'
$ws.Range("E22").Value = '  +1.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.51'
$ws.Range("E23").Value = '  +2.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("E25").Value = '  -1.06%  '

$ws.Range("E26").Value = '  +1.06%  '

$ws.Range("E27").Value = '  +16.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.80'
$ws.Range("E28").Value = '  +5.50%  '

$ws.Range("E29").Value = '  +1.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.90'
$ws.Range("E30").Value = '  +1.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.09'
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.126'
$ws.Range("E32").Value = '  -4.50%  '

$ws.Range("E33").Value = '  +3.99%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0757'
$ws.Range("E34").Value = '  +7.19%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("E35").Value = '  +7.03%  '

$ws.Range("E36").Value = '  +5.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.73'
$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("E38").Value = '  -1.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.37'
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0280'
$ws.Range("E40").Value = '  +5.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.219'
$ws.Range("E41").Value = '  +20.82%  '

$ws.Range("E42").Value = '  +11.35%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.09'
$ws.Range("E43").Value = '  -1.15%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.08'
$ws.Range("E44").Value = '  +2.02%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.71'
$ws.Range("E46").Value = '  +6.00%  '

$ws.Range("E47").Value = '  +9.65%  '

$ws.Range("E48").Value = '  +2.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.86'
$ws.Range("E49").Value = '  +2.36%  '

$ws.Range("E50").Value = '  +0.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.70'
$ws.Range("E51").Value = '  +7.06%  '
